# "Verify Product Version" is the 3rd sheet (sheetId=3 / rId3) - the one whose
# C2 value is being updated from the placeholder number 123 to the real
# product-version string, and whose column widths / selection were tweaked
# after the edit (matches the reviewed xml diff for xl/worksheets/sheet3.xml).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Verify Product Version")

# Replace the placeholder numeric value with the actual product version text
# -> becomes a new shared-string entry ("actiTIME 2017.4").
$ws.Range("C2").Value = "actiTIME 2017.4"

# Widen columns B and C so the longer text is readable.
$ws.Columns.Item(2).ColumnWidth = 14.5
$ws.Columns.Item(3).ColumnWidth = 29.8333333333333

# Leave the cursor/selection on E6, like in the saved workbook.
$ws.Range("E6").Select()
